# Update the dSF column (F) values for rows 2-8 to reflect the repulled
# data / mean calculation referenced in the commit message.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F2").Value = 1
$ws.Range("F3").Value = -1
$ws.Range("F4").Value = -1
$ws.Range("F5").Value = -1
$ws.Range("F6").Value = 11
$ws.Range("F7").Value = 2
$ws.Range("F8").Value = 2
